$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts existing B:G to C:H)
$ws.Columns("B").Insert()

# New column B header is a date like the existing B/C/D headers - copy their number format
$ws.Range("B1").NumberFormat = $ws.Range("C1").NumberFormat

# Fill in the new column B with the new day's data
$ws.Range("B1").Value = 43906
$ws.Range("B2").Value = 14649
$ws.Range("B3").Value = 3522
$ws.Range("B4").Value = 2473
$ws.Range("B5").Value = 1242
$ws.Range("B6").Value = 1516
$ws.Range("B7").Value = 866
$ws.Range("B8").Value = 667
$ws.Range("B9").Value = 523
$ws.Range("B10").Value = 386
$ws.Range("B11").Value = 378
$ws.Range("B12").Value = 400
$ws.Range("B13").Value = 230
$ws.Range("B14").Value = 241
$ws.Range("B15").Value = 213
$ws.Range("B16").Value = 164
$ws.Range("B17").Value = 176
$ws.Range("B18").Value = 107
$ws.Range("B19").Value = 89
$ws.Range("B20").Value = 105
$ws.Range("B21").Value = 21
$ws.Range("B22").Value = 12

# Match the new column widths (A and B both 18.5) and selection/view shown in the saved file
$ws.Columns("A").ColumnWidth = 17.666666666666668
$ws.Columns("B").ColumnWidth = 17.666666666666668

$ws.Range("G14").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 6
